# Update LR-pair (Pgf-Nrp1) sheet with refreshed TPM-derived values.
# The source row set also grows from 12 data rows (clusters ECs/FAPs/MuSCs)
# to 16 data rows, adding the "Resolving-Mac" sending cluster block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$data = @(
    @("ECs","Pgf","Nrp1","ECs",3,1,3.529548333333333,10.588645,0.6379922832219154,0.6379922832219154,3,1,133.7780026666667,401.334008,0.50863533211804,0.5086353321180399,472.1759263487955,4249.58333713916,0.3245054168653256,0.3245054168653255),
    @("ECs","Pgf","Nrp1","FAPs",3,1,3.529548333333333,10.588645,0.6379922832219154,0.6379922832219154,3,1,52.441971,157.325913,0.1993888292903622,0.1993888292903622,185.096471339765,1665.868242057885,0.1272085344479029,0.1272085344479029),
    @("ECs","Pgf","Nrp1","MuSCs",3,1,3.529548333333333,10.588645,0.6379922832219154,0.6379922832219154,3,1,21.197691,63.593073,0.08059542216956049,0.08059542216956046,74.818274939565,673.364474456085,0.05141925740719207,0.05141925740719205),
    @("ECs","Pgf","Nrp1","Resolving-Mac",3,1,3.529548333333333,10.588645,0.6379922832219154,0.6379922832219154,3,1,55.59592133333333,166.787764,0.2113804164220374,0.2113804164220373,196.2284914821977,1766.05642333978,0.1348590745014949,0.1348590745014948),
    @("FAPs","Pgf","Nrp1","ECs",2,0.6666666666666666,0.6384863333333334,1.915459,0.1154111844176443,0.1154111844176443,3,1,133.7780026666667,401.334008,0.50863533211804,0.5086353321180399,85.4154264032969,768.738837629672,0.05870220611640486,0.05870220611640484),
    @("FAPs","Pgf","Nrp1","FAPs",2,0.6666666666666666,0.6384863333333334,1.915459,0.1154111844176443,0.1154111844176443,3,1,52.441971,157.325913,0.1993888292903622,0.1993888292903622,33.483481776563,301.351335989067,0.02301170094804818,0.02301170094804818),
    @("FAPs","Pgf","Nrp1","MuSCs",2,0.6666666666666666,0.6384863333333334,1.915459,0.1154111844176443,0.1154111844176443,3,1,21.197691,63.593073,0.08059542216956049,0.08059542216956046,13.534436001723,121.809924015507,0.009301613131229042,0.009301613131229038),
    @("FAPs","Pgf","Nrp1","Resolving-Mac",2,0.6666666666666666,0.6384863333333334,1.915459,0.1154111844176443,0.1154111844176443,3,1,55.59592133333333,166.787764,0.2113804164220374,0.2113804164220373,35.49723596040845,319.4751236436759,0.0243956642219622,0.02439566422196219),
    @("MuSCs","Pgf","Nrp1","ECs",3,1,1.109135333333333,3.327406,0.2004845144158011,0.2004845144158011,3,1,133.7780026666667,401.334008,0.50863533211804,0.5086353321180399,148.3779095803609,1335.401186223248,0.101973507574405,0.101973507574405),
    @("MuSCs","Pgf","Nrp1","FAPs",3,1,1.109135333333333,3.327406,0.2004845144158011,0.2004845144158011,3,1,52.441971,157.325913,0.1993888292903622,0.1993888292903622,58.165242985742,523.487186871678,0.03997437262021334,0.03997437262021333),
    @("MuSCs","Pgf","Nrp1","MuSCs",3,1,1.109135333333333,3.327406,0.2004845144158011,0.2004845144158011,3,1,21.197691,63.593073,0.08059542216956049,0.08059542216956046,23.511108073182,211.599972658638,0.01615813407780083,0.01615813407780082),
    @("MuSCs","Pgf","Nrp1","Resolving-Mac",3,1,1.109135333333333,3.327406,0.2004845144158011,0.2004845144158011,3,1,55.59592133333333,166.787764,0.2113804164220374,0.2113804164220373,61.66340074002044,554.9706066601839,0.042378500143382,0.04237850014338198),
    @("Resolving-Mac","Pgf","Nrp1","ECs",1,0.3333333333333333,0.2551043333333333,0.765313,0.04611201794463916,0.04611201794463916,3,1,133.7780026666667,401.334008,0.50863533211804,0.5086353321180399,34.12734818494489,307.146133664504,0.02345420156190456,0.02345420156190456),
    @("Resolving-Mac","Pgf","Nrp1","FAPs",1,0.3333333333333333,0.2551043333333333,0.765313,0.04611201794463916,0.04611201794463916,3,1,52.441971,157.325913,0.1993888292903622,0.1993888292903622,13.378174050641,120.403566455769,0.009194221274197776,0.009194221274197776),
    @("Resolving-Mac","Pgf","Nrp1","MuSCs",1,0.3333333333333333,0.2551043333333333,0.765313,0.04611201794463916,0.04611201794463916,3,1,21.197691,63.593073,0.08059542216956049,0.08059542216956046,5.407622830761,48.668605476849,0.003716417553338542,0.003716417553338541),
    @("Resolving-Mac","Pgf","Nrp1","Resolving-Mac",1,0.3333333333333333,0.2551043333333333,0.765313,0.04611201794463916,0.04611201794463916,3,1,55.59592133333333,166.787764,0.2113804164220374,0.2113804164220373,14.18276044779244,127.644844030132,0.009747177555198284,0.009747177555198284)
)

$rowIndex = 2
foreach ($row in $data) {
    for ($colIndex = 0; $colIndex -lt $columns.Length; $colIndex++) {
        $cellRef = "$($columns[$colIndex])$rowIndex"
        $ws.Range($cellRef).Value = $row[$colIndex]
    }
    $rowIndex++
}
